# Update NG_PROD_SUM_DC_NUS_MMCF_A.xlsx with refreshed EIA data (2019 annual
# release), per commit "add updated NG data from EIA".

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# "Contents" sheet — release/date metadata housekeeping.
# ----------------------------------------------------------------------
$wsContents = $wb.Worksheets.Item("Contents")

# Latest-data year shown in the summary strip (was 2017).
$wsContents.Range("F7").Value = 2019

# These three cells hold date-looking strings that must stay TEXT (not get
# auto-converted to date serials) — use a leading apostrophe so Excel keeps
# them as literal text, same as the source file.
$wsContents.Range("G7").Value = "'6/30/1900"
$wsContents.Range("C9").Value = "'7/31/2020"
$wsContents.Range("C10").Value = "'8/31/2020"
$wsContents.Range("F16").Value = "'7/29/2020 8:36:57 PM"

# ----------------------------------------------------------------------
# "Data 1" sheet — refreshed annual series through 6/30/2019, plus a
# revision to the two most recent prior years.
# ----------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data 1")

# Row 120 (6/30/2016) — revised figures.
$wsData.Cells.Item(120, 2).Value = 32591578
$wsData.Cells.Item(120, 3).Value = 7287858
$wsData.Cells.Item(120, 4).Value = 6385120
$wsData.Cells.Item(120, 5).Value = 17847539
$wsData.Cells.Item(120, 6).Value = 1071062
$wsData.Cells.Item(120, 7).Value = 3548106
$wsData.Cells.Item(120, 8).Value = 230410
$wsData.Cells.Item(120, 9).Value = 413013
$wsData.Cells.Item(120, 10).Value = 28400049
$wsData.Cells.Item(120, 11).Value = 1807934
$wsData.Cells.Item(120, 12).Value = 26592115

# Row 121 (6/30/2017) — revised/filled-in figures.
$wsData.Cells.Item(121, 2).Value = 33292113
$wsData.Cells.Item(121, 3).Value = 6161420
$wsData.Cells.Item(121, 4).Value = 6217438
$wsData.Cells.Item(121, 5).Value = 19927602
$wsData.Cells.Item(121, 6).Value = 985653
$wsData.Cells.Item(121, 9).Value = 264582
$wsData.Cells.Item(121, 10).Value = 29203550
$wsData.Cells.Item(121, 11).Value = 1897242
$wsData.Cells.Item(121, 12).Value = 27306308

# Row 122 (6/30/2018) — new row.
$wsData.Cells.Item(122, 1).Value = 43281
$wsData.Cells.Item(122, 2).Value = 37129374
$wsData.Cells.Item(122, 3).Value = 6350001
$wsData.Cells.Item(122, 4).Value = 6275713
$wsData.Cells.Item(122, 5).Value = 23550471
$wsData.Cells.Item(122, 6).Value = 953189
$wsData.Cells.Item(122, 7).Value = 3584274
$wsData.Cells.Item(122, 8).Value = 468347
$wsData.Cells.Item(122, 9).Value = 253459
$wsData.Cells.Item(122, 10).Value = 32823295
$wsData.Cells.Item(122, 11).Value = 2234593
$wsData.Cells.Item(122, 12).Value = 30588702

# Row 123 (6/30/2019) — new row.
$wsData.Cells.Item(123, 1).NumberFormat = "yyyy"
$wsData.Cells.Item(123, 1).Value = 43646
$wsData.Cells.Item(123, 2).Value = 40704488
$wsData.Cells.Item(123, 10).Value = 36197056
$wsData.Cells.Item(123, 11).Value = 2540010
$wsData.Cells.Item(123, 12).Value = 33657046

# Row 124 — new trailing blank row (matches the style of the date column
# used throughout, so column A keeps its "yyyy" number format).
$wsData.Cells.Item(124, 1).NumberFormat = "yyyy"

# The refreshed workbook re-opens on the "Contents" tab rather than "Data 1".
$wsContents.Activate()
